# Weekly update: insert a new price record as the new row 295 in the
# "Acelga" (Hortaliza) table, pushing the existing rows 295:363 down to
# 296:364.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 295; everything that was
# in rows 295:363 shifts down to 296:364.
$ws.Rows.Item(295).Insert()

# Populate the new row 295 with the latest weekly observation. It mirrors
# the row that used to be at 295 (same market/region/category/etc.) but
# with an updated date and volume.
$ws.Range("A295").Value = 4
$ws.Range("B295").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C295").Value = "Los Lagos"
$ws.Range("D295").Value = 45275
$ws.Range("E295").Value = 10
$ws.Range("F295").Value = 100112009
$ws.Range("G295").Value = "Acelga"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 75
$ws.Range("K295").Value = 10000
$ws.Range("L295").Value = 10000
$ws.Range("M295").Value = 10000
$ws.Range("N295").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O295").Value = "Región de La Araucanía"
$ws.Range("P295").Value = 833
$ws.Range("Q295").Value = 12
$ws.Range("R295").Value = "Hortaliza"
